$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.925.97"
$ws.Range("E2").Value = '  -2.08%  '
$ws.Range("D3").Value = "'2.463.84"
$ws.Range("E3").Value = '  -2.53%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = "'517.44"
$ws.Range("E5").Value = '  -3.71%  '
$ws.Range("D6").Value = "'130.61"
$ws.Range("E6").Value = '  -4.91%  '
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  -2.18%  '
$ws.Range("E9").Value = '  -2.44%  '
$ws.Range("E10").Value = '  -0.61%  '
$ws.Range("D11").Value = "'5.33"
$ws.Range("E11").Value = '  +0.10%  '
$ws.Range("E12").Value = '  -2.13%  '
$ws.Range("D13").Value = "'2.900.97"
$ws.Range("E13").Value = '  -2.41%  '
$ws.Range("D14").Value = "'57.852.50"
$ws.Range("E14").Value = '  -2.09%  '
$ws.Range("D15").Value = "'22.25"
$ws.Range("E15").Value = '  -3.62%  '
$ws.Range("E16").Value = '  -2.60%  '
$ws.Range("D17").Value = "'2.463.27"
$ws.Range("E17").Value = '  -2.78%  '
$ws.Range("D18").Value = "'10.74"
$ws.Range("E18").Value = '  -3.89%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = "'319.93"
$ws.Range("E19").Value = '  -1.26%  '
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").Value = "'4.17"
$ws.Range("E20").Value = '  -2.99%  '
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("D22").Value = "'5.72"
$ws.Range("E22").Value = '  -3.94%  '
$ws.Range("D23").Value = "'64.01"
$ws.Range("E23").Value = '  -2.96%  '
$ws.Range("E24").Value = '  -3.37%  '
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("E26").Value = '  -3.31%  '
$ws.Range("D27").Value = "'7.30"
$ws.Range("E27").Value = '  -3.38%  '
$ws.Range("D28").Value = "'0.0₃0750"
$ws.Range("E28").Value = '  -3.39%  '
$ws.Range("E29").Value = '  -4.97%  '
$ws.Range("D30").Value = "'165.60"
$ws.Range("E30").Value = '  -0.72%  '
$ws.Range("D31").Value = "'6.29"
$ws.Range("E31").Value = '  -6.40%  '
$ws.Range("E32").Value = '  -2.81%  '
$ws.Range("D33").Value = "'0.998"
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = '  +0.05%  '
$ws.Range("D35").Value = "'18.02"
$ws.Range("E35").Value = '  -2.30%  '
$ws.Range("E37").Value = '  -3.59%  '
$ws.Range("E38").Value = '  -5.00%  '
$ws.Range("D39").Value = "'0.787"
$ws.Range("E39").Value = '  -3.70%  '
$ws.Range("E40").Value = '  -5.02%  '
$ws.Range("D41").Value = "'271.61"
$ws.Range("E41").Value = '  -4.87%  '
$ws.Range("E42").Value = '  -3.14%  '
$ws.Range("E43").Value = '  -2.69%  '
$ws.Range("D44").Value = "'125.70"
$ws.Range("E44").Value = '  -5.13%  '
$ws.Range("E45").Value = '  -2.36%  '
$ws.Range("D46").Value = "'0.0487"
$ws.Range("E46").Value = '  -4.57%  '
$ws.Range("E47").Value = '  -3.75%  '
$ws.Range("D48").Value = "'17.00"
$ws.Range("E48").Value = '  -2.39%  '
$ws.Range("D49").Value = "'1.728.20"
$ws.Range("E49").Value = '  -2.24%  '
$ws.Range("D50").Value = "'0.976"
$ws.Range("E50").Value = '  -1.45%  '
$ws.Range("E51").Value = '  -1.28%  '
